$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(10616,10616,10616,10602,10602,10192,10192,9927,9927,8887,8887,8887,8842,8842,8421,8421,8421,8421,8421,8421,8257,8257,8257,8257,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,8204,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7678,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653,7653)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

